# Apply weekly re-ordering of Fecha (D), Volumen (M), Precio minimo (N),
# Precio maximo (O), Precio promedio ponderado (P) and Precio $/Kg (S)
# across rows 2-8. Other columns remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get reshuffled per row.
$cols = @("D", "M", "N", "O", "P", "S")

# Capture current ("before") values for the affected columns/rows first,
# since rows will be overwritten using values from other rows.
$rows = 2..8
$before = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# Destination row -> source row mapping (values move from source to destination).
$mapping = @{
    2 = 6
    3 = 7
    4 = 5
    5 = 2
    6 = 8
    7 = 3
    8 = 4
}

foreach ($dst in $mapping.Keys) {
    $src = $mapping[$dst]
    foreach ($c in $cols) {
        $ws.Range("$c$dst").Value2 = $before[$src][$c]
    }
}
